$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source data added a new account row ("005081833", "PEDRO", 100000)
# right above the existing "004526450" (MSD) row, i.e. at sheet row 5
# (row 1 is the header "Conta"/"Nome"/"Saldo"). Insert a fresh row there
# so every row below shifts down by one, matching the diff.
$ws.Rows.Item(5).Insert()

# Account numbers are stored as text (leading zeros must be preserved),
# so force text formatting just long enough to type the value, then
# clear the format back off again so the cell ends up with the same
# (default/no) style as its sibling cells.
$ws.Range("A5").NumberFormat = "@"
$ws.Range("A5").Value = "005081833"
$ws.Range("A5").ClearFormats()

$ws.Range("B5").Value = "PEDRO"
$ws.Range("C5").Value = 100000
